# Updated jxls 2 report in demo db
#
# Re-targets the sample jxls2 template from the old jxls1-style expression
# syntax (${results.X}, ${order_date}, ${city}, ${item}) to the jxls2
# syntax (${row.X}, ${order_date.value}, ${cityId...}, ${item...}) and adds
# the jx:area / jx:each control comments that drive the jxls2 template
# processor, plus a couple of layout tweaks (column widths, date format,
# selected cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Parameters block (rows 4-6)
# ---------------------------------------------------------------------
$ws.Range("A4").Value = '${order_date.value}'
$ws.Range("A4").NumberFormat = "yyyy\-mm\-dd;@"

$ws.Range("A5").Value = '${cityId.parameter.label}'
$ws.Range("B5").Value = '${cityId.value}'

$ws.Range("A6").Value = '${item.nameAndDisplayValues}'

# ---------------------------------------------------------------------
# Table header (row 9) - columns reordered: Order ID, City, Item, Date, Volume
# ---------------------------------------------------------------------
$ws.Range("A9").Value = "Order ID"
$ws.Range("B9").Value = "City"
$ws.Range("C9").Value = "Item"
$ws.Range("D9").Value = "Date"
$ws.Range("E9").Value = "Volume"

# ---------------------------------------------------------------------
# Table data row (row 10) - jxls2 "row" loop variable instead of "results"
# ---------------------------------------------------------------------
$ws.Range("A10").Value = '${row.ORDER_ID}'
$ws.Range("B10").Value = '${row.CITY_NAME}'
$ws.Range("C10").Value = '${row.ITEM_NAME}'
$ws.Range("D10").Value = '${row.ORDER_DATE}'
$ws.Range("E10").Value = '${row.VOLUME}'

# ---------------------------------------------------------------------
# Column widths for the newly meaningful A:C columns
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 23.786830357142858
$ws.Columns.Item(2).ColumnWidth = 16.501116071428573
$ws.Columns.Item(3).ColumnWidth = 16.358258928571427

# ---------------------------------------------------------------------
# jxls directives as cell comments
# ---------------------------------------------------------------------
$ws.Range("A1").AddComment("Author:" + [char]10 + 'jx:area(lastCell="E10")') | Out-Null
$ws.Range("A10").AddComment("Author:" + [char]10 + 'jx:each(items="results" var="row" lastCell="E10")') | Out-Null

# ---------------------------------------------------------------------
# Selected cell
# ---------------------------------------------------------------------
$ws.Range("A7").Select() | Out-Null
